$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
# manual shift: M17 = old M16, M18 = old M17, M16 = new value
$ws.Range("M18").Value = $ws.Range("M17").Value
$ws.Range("M17").Value = $ws.Range("M16").Value
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"
Write-Output "done"
